$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '68.115.44'
$ws.Range("E2").Value = '  +1.25%  '

$ws.Range("D3").Value = '3.273.13'
$ws.Range("E3").Value = '  +0.97%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '185.34'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.78%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("E9").Value = '  +4.87%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.12%  '

$ws.Range("E11").Value = '  +1.13%  '

$ws.Range("D12").Value = '3.841.49'
$ws.Range("E12").Value = '  +0.99%  '

$ws.Range("E13").Value = '  +0.30%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.59'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.53%  '

$ws.Range("D15").Value = '68.113.39'
$ws.Range("E15").Value = '  +1.36%  '

$ws.Range("D17").Value = '3.274.82'
$ws.Range("E17").Value = '  +0.87%  '

$ws.Range("E18").Value = '  +0.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.63%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '381.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.74%  '

$ws.Range("E21").Value = '  +2.26%  '

$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.38'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.45%  '

$ws.Range("E24").Value = '  +2.76%  '

$ws.Range("E25").Value = '  +1.16%  '

$ws.Range("E26").Value = '  +6.62%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.78'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.85'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.41%  '

$ws.Range("E30").Value = '  +1.30%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.93'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.92%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.19'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.68%  '

$ws.Range("E33").Value = '  +0.04%  '

$ws.Range("E34").Value = '  +0.70%  '

$ws.Range("E35").Value = '  +3.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.22'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.83%  '

$ws.Range("E37").Value = '  +0.59%  '

$ws.Range("E38").Value = '  -1.96%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.75'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.61'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.72%  '

$ws.Range("E41").Value = '  +1.50%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.60'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.80%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.21%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.34'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.21%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0690'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.89%  '

$ws.Range("D46").Value = '2.635.93'
$ws.Range("E46").Value = '  -4.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '342.10'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '32.27'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.28%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.05%  '

$ws.Range("E51").Value = '  -0.19%  '
